# Fix 3-Year Summary category alignment: Azure DevOps Enterprise Platform
# - Row 5 label changes from "Support & Maintenance" to "Hardware"
# - Row 6 becomes a new "Support & Maintenance" row (SUMIF formulas, same
#   style the old TOTAL row already had)
# - Row 7 is a new "Professional Services" row (default styling)
# - Row 8 is the TOTAL row, relocated down, SUM ranges widened to B3:B7
# - Row 9 is left as a trailing blank row, matching the sheet's dimension

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3-Year Summary")

# --- Row 5: rename "Support & Maintenance" -> "Hardware" (style/formulas unchanged) ---
$ws.Range("A5").Value = "Hardware"

# --- Row 6: was the TOTAL row; reuse it for the (relocated) "Support & Maintenance" ---
# Row 6 already carries the s="51"/s="53" styling the new content needs, so only the
# cell contents are rewritten.
$ws.Range("A6").Value = "Support & Maintenance"
$ws.Range("B6").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A6,'Infrastructure Costs'!`$G:`$G)"
$ws.Range("C6").Formula = "=SUMIF(Credits!`$A:`$A,A6,Credits!`$C:`$C)"
$ws.Range("D6").Formula = "=B6+C6"
$ws.Range("E6").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A6,'Infrastructure Costs'!`$H:`$H)"
$ws.Range("F6").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A6,'Infrastructure Costs'!`$I:`$I)"
$ws.Range("G6").Formula = "=D6+E6+F6"

# --- Row 7: new "Professional Services" row (plain default styling) ---
$ws.Range("A7").Value = "Professional Services"
$ws.Range("B7").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A7,'Infrastructure Costs'!`$G:`$G)"
$ws.Range("C7").Formula = "=SUMIF(Credits!`$A:`$A,A7,Credits!`$C:`$C)"
$ws.Range("D7").Formula = "=B7+C7"
$ws.Range("E7").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A7,'Infrastructure Costs'!`$H:`$H)"
$ws.Range("F7").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A7,'Infrastructure Costs'!`$I:`$I)"
$ws.Range("G7").Formula = "=D7+E7+F7"

# --- Row 8: TOTAL row, moved down one row; widen SUM ranges to include row 7 ---
$ws.Range("A8").Value = "TOTAL"
$ws.Range("B8").Formula = "=SUM(B3:B7)"
$ws.Range("C8").Formula = "=SUM(C3:C7)"
$ws.Range("D8").Formula = "=SUM(D3:D7)"
$ws.Range("E8").Formula = "=SUM(E3:E7)"
$ws.Range("F8").Formula = "=SUM(F3:F7)"
$ws.Range("G8").Formula = "=SUM(G3:G7)"
# SUM(...) over the currency-formatted rows above auto-inherits their number
# format; the source row had no explicit style, so put row 8 back to default.
$ws.Range("A8:G8").Style = "Normal"

# --- Row 9: trailing blank row so the sheet's dimension extends to G9 ---
$ws.Range("A9").Style = "Normal"

Write-Output "3-Year Summary rows 5-9 updated"
